$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.644.59'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.09%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.635.09'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.39%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.83%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.535'
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.634.46'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.35%  '

$ws.Range("E10").Value = '  -2.67%  '

$ws.Range("E11").Value = '  +1.34%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.365'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.77%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.23'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.43%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.70'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.64%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.117.62'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.56%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000182'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.08%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.600.65'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.12%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.618.54'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.98%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.04'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.51%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '357.99'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.34'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.70'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.19%  '

$ws.Range("E24").Value = '  -3.83%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.38'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.40%  '

$ws.Range("E26").Value = '  +0.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '70.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.08%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.770.05'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.79%  '

$ws.Range("E29").Value = '  +0.34%  '

$ws.Range("E30").Value = '  -1.17%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '550.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.30%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.96'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.62%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.36'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.90'
$ws.Range("D34").Style = "Normal"

$ws.Range("E35").Value = '  +5.45%  '

$ws.Range("E36").Value = '  +0.10%  '

$ws.Range("E37").Value = '  -2.73%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '158.33'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.03%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.09'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.367'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.41%  '

$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.82'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.89%  '

$ws.Range("B42").Value = 'WhiteBITCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '18.29'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.02%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.23'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.20%  '

$ws.Range("E44").Value = '  +0.15%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.44'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.91%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₆0302'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.04%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '153.13'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.39%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.583'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.18%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.33%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.71'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.85%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0773'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.71%  '
